# Generate Report for Handoff
# Rewrites the UUID-named source file references from the old GUID
# (649e5cf0-066b-432c-bf2b-c3a9f3c991ff) to the new one
# (3ff162d2-12b8-4ca1-af4e-ab11fcddb641), refreshes the handoff xliff
# filenames/timestamps, and clears the (now stale) handback target/file/date
# info while flagging "Has metadata" = True on both locale sheets.

$wb = $excel.ActiveWorkbook

$oldGuid = "649e5cf0-066b-432c-bf2b-c3a9f3c991ff"
$newGuid = "3ff162d2-12b8-4ca1-af4e-ab11fcddb641"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newGuid.md"

foreach ($hl in $wsOverview.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$B$2') {
        $hl.TextToDisplay = "e2e\$newGuid.md"
    }
}
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"

# "Latest HO Xliff Generate Date" mirrors the de-de handoff timestamp below.
$wsOverview.Range("G2").Value = "2016-11-09 06:41:15"

$wsOverview.Columns.Item(1).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid.md"
    }
}
$wsZh.Range("A2").Value = "$newGuid.md"

$wsZh.Range("G2").Value = "$newGuid.8ca6a6e7f5ad85b6f220aa3c5490c4f62cfe9d23.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-11-09 06:41:01"

# Remove the now-stale "Latest Target File" hyperlink at I2, keeping the
# one at A2 intact.
foreach ($hl in $wsZh.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2') {
        $hl.Delete()
    }
}
$wsZh.Range("I2").Value = "'"
$wsZh.Range("I2").Style = "Normal"

$wsZh.Range("J2").Value = "'"
$wsZh.Range("J2").Style = "Normal"

$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

$wsZh.Range("O2").Value = "'True"
$wsZh.Range("O2").Style = "Normal"

$wsZh.Columns.Item(1).ColumnWidth = 39.17
$wsZh.Columns.Item(9).ColumnWidth = 17.8
$wsZh.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$A$2') {
        $hl.TextToDisplay = "$newGuid.md"
    }
}
$wsDe.Range("A2").Value = "$newGuid.md"

$wsDe.Range("G2").Value = "$newGuid.8ca6a6e7f5ad85b6f220aa3c5490c4f62cfe9d23.de-de.xlf"
$wsDe.Range("H2").Value = "2016-11-09 06:41:15"

foreach ($hl in $wsDe.Hyperlinks) {
    $addr = $hl.Range.Address()
    if ($addr -eq '$I$2') {
        $hl.Delete()
    }
}
$wsDe.Range("I2").Value = "'"
$wsDe.Range("I2").Style = "Normal"

$wsDe.Range("J2").Value = "'"
$wsDe.Range("J2").Style = "Normal"

$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$wsDe.Range("O2").Value = "'True"
$wsDe.Range("O2").Style = "Normal"

$wsDe.Columns.Item(1).ColumnWidth = 39.17
$wsDe.Columns.Item(9).ColumnWidth = 17.8
$wsDe.Columns.Item(10).ColumnWidth = 20.8
